# implemented new partition method - link process method
#
# Recreates the "link process method" partitioning results: the previously
# all-zero emission rows are cleared out entirely (leaving only the region
# label in column A), and the non-zero rows are updated with the values
# produced by the new partition/link method. Also restores the view to the
# default top-left cell instead of the old B2:D2 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- window / view bookkeeping -------------------------------------------
$win = $excel.ActiveWindow
$win.Width = 19905
$win.Height = 6810

# --- remove rows that are now fully empty (were all zeros) ---------------
$ws.Range("B2:D2").ClearContents()
$ws.Range("B3:D3").ClearContents()
$ws.Range("B4:D4").ClearContents()
$ws.Range("B5:D5").ClearContents()
$ws.Range("B6:D6").ClearContents()
$ws.Range("B8:D8").ClearContents()
$ws.Range("B9:D9").ClearContents()
$ws.Range("B10:D10").ClearContents()
$ws.Range("B13:D13").ClearContents()
$ws.Range("B15:D15").ClearContents()
$ws.Range("B16:D16").ClearContents()

# --- updated emission-factor values from the new partition method --------
$ws.Range("B7").Value = 2057.5853342036689
$ws.Range("C7").Value = 73187.27039676116
$ws.Range("D7").Value = 1637.0359030066611

$ws.Range("B11").Value = 741.60154071836155
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 60058.615613782968

$ws.Range("B12").Value = 1021.0115996368736
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 60042.538077040212

$ws.Range("B14").Value = 1017.7456914109041
$ws.Range("C14").Value = 22007.769192199543
$ws.Range("D14").Value = 32074.47201424107

$ws.Range("B17").Value = 3022.1657283353525
$ws.Range("C17").Value = 106207.13997045124
$ws.Range("D17").Value = 2521.7475685541099

$ws.Range("B18").Value = 5968.872656911175
$ws.Range("C18").Value = 108685.87564136037
$ws.Range("D18").Value = 5733.023327330091

$ws.Range("B19").Value = 764.37264664711483
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 61768.08185027646

$ws.Range("B20").Value = 6463.9516945768892
$ws.Range("C20").Value = 212192.50511015704
$ws.Range("D20").Value = 7660.6230599706832

$ws.Range("B21").Value = 87.901497846762595
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 76658.826509068138

$ws.Range("B22").Value = 743.12690180320158
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 66158.105600168492

$ws.Range("B23").Value = 808.95857475450384
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 60665.239349189083

$ws.Range("B24").Value = 824.15457453162492
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 63007.538830683436

$ws.Range("B25").Value = 784.45259281914196
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 60742.167990994021

$ws.Range("B26").Value = 867.05002640084581
$ws.Range("C26").Value = 27630.948019681993
$ws.Range("D26").Value = 50453.645874829395

$ws.Range("B27").Value = 4200.9117495755672
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 59868.309374172648

$ws.Range("B28").Value = 6903.2641583531558
$ws.Range("C28").Value = 259860.99340165633
$ws.Range("D28").Value = 6488.0581168295048

# --- reset the selection to the top-left cell (was B2:D2) ----------------
[void]$ws.Range("A1").Select()
